$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: apply the textual changes (content level). Word's Find/Replace
# re-normalises (merges) the runs of whichever paragraph it touches, so
# run-level splitting is re-established afterwards, in Step 2.
# ---------------------------------------------------------------------

# "... goes through a variety of physical ..." -> "... goes through different physical ..."
$d.Content.Find.Execute("goes through a variety of physical", $true, $true, $false, $false, $false, $true, 1, $false, "goes through different physical", 2) | Out-Null

# "... these “target” sounds and separate them from this ..." -> "... these “target” sound and separate it from this ..."
$d.Content.Find.Execute("sounds and separate them from", $true, $true, $false, $false, $false, $true, 1, $false, "sound and separate it from", 2) | Out-Null

# "(so called “auditory objects”)" -> "(“auditory objects”)"
$d.Content.Find.Execute("so called “auditory objects", $true, $true, $false, $false, $false, $true, 1, $false, "“auditory objects", 2) | Out-Null

# "... underlying physical and biological processes." -> "... underlying physics and biology."
$d.Content.Find.Execute("underlying physical and biological processes.", $true, $true, $false, $false, $false, $true, 1, $false, "underlying physics and biology.", 2) | Out-Null

# ---------------------------------------------------------------------
# Step 2: re-split the runs so the affected paragraphs end up with the
# same run boundaries Word itself would leave after this kind of
# in-place retyping (new/changed words get their own <w:r>, the
# untouched text around them keeps its own <w:r> too). Toggling a
# formatting property and back forces the engine to keep the touched
# sub-range from re-merging into its identically-formatted neighbours.
# ---------------------------------------------------------------------

function Split-Boundaries($paragraphRange, [string[]]$chunks) {
    $base = $paragraphRange.Start
    $pos = 0
    for ($i = 0; $i -lt $chunks.Length; $i++) {
        $chunkStart = $base + $pos
        $pos += $chunks[$i].Length
        $chunkEnd = $base + $pos
        # Only interior chunks need a forced split on both sides; the
        # very first/last chunk is already bounded by the paragraph
        # edges, but toggling them too is harmless and keeps the logic
        # simple/uniform.
        $r = $d.Range($chunkStart, $chunkEnd)
        $r.Font.Bold = 1
        $r.Font.Bold = 0
    }
}

# Paragraph: "Imagine a party. ... ASA."
$p1 = $d.Paragraphs(3).Range
Split-Boundaries $p1 @(
    "Imagine a party. You can hear a wide variety of sounds: music",
    " ",
    "i",
    "n the background",
    ", ",
    "conversations",
    " between people",
    ", noises",
    " of somebody coughing,",
    " maybe even a dog barking outside",
    "… These sounds form a ",
    "single",
    " stream, which comes to ",
    "y",
    "our ears in a form of a sound wave and then goes through ",
    "different",
    " physical, biological and psychoacoustical processes to the brain. Despite all these sounds ",
    "from different sources are ",
    "mixed on the way to your ears, your brain can segregate one (or several) of them. You can focus your hearing on these “target” sound and separate ",
    "it",
    " from this complex mixture, ",
    "leaving other sounds in the background. ",
    "This phenomenon",
    " has been described as a “cocktail party ",
    "effect",
    "”, and the process of ",
    "integrating separate sounds into meaningful streams (",
    "“auditory objects”",
    ") --",
    " auditory scene analysis",
    ", or ",
    "ASA."
)

# Paragraph: "In machine perception ... underlying physics and biology."
$p2 = $d.Paragraphs(4).Range
Split-Boundaries $p2 @(
    "In machine perception",
    " ---",
    " specifically in machine hearing",
    " --- the related concept is referred to as Computational ASA (CASA)",
    " and ",
    "is tightly connected ",
    "to the ",
    "fields",
    " of sound recognition and digital signal processing.",
    " ",
    "The main objective for this thesis is to describe ",
    "its",
    " principles and goals",
    ", ",
    "along with existing applications and ",
    "approaches. ",
    "Another objective would be to practically apply the theoretical knowledge and implement a simple CASA system to separate monophonic music from noise.",
    " ",
    "But firstly",
    ", since this thesis is made for an IT-oriented ",
    "audience, ",
    "it is needed to ",
    "make a brief introduction",
    " to the underlying physic",
    "s",
    " and biolog",
    "y",
    "."
)
